$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($Row, $D, $J, $K, $L, $M, $N, $O, $P, $Q) {
    $ws.Cells.Item($Row, 1).Value = 10
    $ws.Cells.Item($Row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($Row, 3).Value = "La Araucanía"
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = 9
    $ws.Cells.Item($Row, 6).Value = 100114002
    $ws.Cells.Item($Row, 7).Value = "Camote"
    $ws.Cells.Item($Row, 8).Value = "Sin especificar"
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# Insert a new record before the current row 158 (new weekly price entry,
# Peru, $/caja 18 kilos), pushing the existing rows 158-217 down to 159-218.
$ws.Rows.Item(158).Insert()
Set-DataRow 158 45120 100 26000 26000 26000 "`$/caja 18 kilos" "Perú" 1444 18

# Insert a second new record before the (now shifted) row 214, pushing the
# remaining rows 214-218 down to 215-219.
$ws.Rows.Item(214).Insert()
Set-DataRow 214 45121 50 26000 26000 26000 "`$/caja 18 kilos" "Perú" 1444 18
